$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.505.66"
$ws.Range("E2").Value = "  +3.71%  "

# Row 3
$ws.Range("D3").Value = "2.421.32"
$ws.Range("E3").Value = "  +2.54%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.82"
$ws.Range("E5").Value = "  +3.86%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.33"
$ws.Range("E6").Value = "  +5.82%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.511"
$ws.Range("E7").Value = "  +1.96%  "

# Row 8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  +5.14%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.29"
$ws.Range("E10").Value = "  +3.56%  "

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.126"
$ws.Range("E11").Value = "  +1.86%  "

# Row 12
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0799"
$ws.Range("E12").Value = "  +2.02%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.94"
$ws.Range("E13").Value = "  +3.51%  "

# Row 14
$ws.Range("E14").Value = "  +3.33%  "

# Row 15
$ws.Range("D15").Value = "2.799.79"
$ws.Range("E15").Value = "  +2.56%  "

# Row 16
$ws.Range("D16").Value = "2.422.19"
$ws.Range("E16").Value = "  +2.34%  "

# Row 17
$ws.Range("E17").Value = "  +5.35%  "

# Row 18
$ws.Range("D18").Value = "44.401.46"
$ws.Range("E18").Value = "  +3.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.46"
$ws.Range("E19").Value = "  +5.52%  "

# Row 20
$ws.Range("E20").Value = "  +2.18%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0922"
$ws.Range("E21").Value = "  +4.56%  "

# Row 22
$ws.Range("E22").Value = "  +1.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.91"
$ws.Range("E23").Value = "  +3.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.29"
$ws.Range("E24").Value = "  +6.17%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +1.70%  "

# Row 26
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.17"
$ws.Range("E27").Value = "  +2.94%  "

# Row 28
$ws.Range("E28").Value = "  -4.34%  "

# Row 29
$ws.Range("E29").Value = "  +3.26%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.41"
$ws.Range("E30").Value = "  +4.78%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.47"
$ws.Range("E31").Value = "  +0.97%  "

# Row 32
$ws.Range("E32").Value = "  +18.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.25"
$ws.Range("E33").Value = "  +11.43%  "

# Row 34
$ws.Range("E34").Value = "  +3.29%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0770"
$ws.Range("E35").Value = "  +8.61%  "

# Row 36
$ws.Range("E36").Value = "  +0.17%  "

# Row 37
$ws.Range("E37").Value = "  +4.69%  "

# Row 38
$ws.Range("E38").Value = "  +2.57%  "

# Row 39
$ws.Range("E39").Value = "  +2.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "121.84"
$ws.Range("E40").Value = "  -3.85%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.22"
$ws.Range("E41").Value = "  -2.53%  "

# Row 42
$ws.Range("E42").Value = "  +1.51%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.01"
$ws.Range("E43").Value = "  -0.98%  "

# Row 44
$ws.Range("E44").Value = "  +4.53%  "

# Row 45
$ws.Range("D45").Value = "1.952.42"
$ws.Range("E45").Value = "  +1.28%  "

# Row 46
$ws.Range("E46").Value = "  +2.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.94"
$ws.Range("E47").Value = "  +8.59%  "

# Row 48
$ws.Range("E48").Value = "  +3.59%  "

# Row 49
$ws.Range("E49").Value = "  +10.31%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.43"
$ws.Range("E50").Value = "  +7.70%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.34"
$ws.Range("E51").Value = "  +4.79%  "
